$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1").Value = "First day - ENEM 2022"
$ws.Range("G1").Value = "Second day - ENEM 2022"

$data = @(
    @(2, 31.65, 36.41),
    @(3, 24.73, 29.48),
    @(4, 42.99, 48.61),
    @(5, 31.91, 37.04),
    @(6, 20.86, 24.73),
    @(7, 20.65, 25.59),
    @(8, 28.8, 33.01),
    @(9, 21.84, 25.44),
    @(10, 17.13, 19.5),
    @(11, 26.66, 30.2),
    @(12, 14.57, 17.57),
    @(13, 23.11, 25.69),
    @(14, 22.79, 26.26),
    @(15, 24.78, 28.26),
    @(16, 19.16, 21.92),
    @(17, 21.3, 24.65),
    @(18, 18.78, 21.94),
    @(19, 19.12, 23),
    @(20, 20.48, 23.94),
    @(21, 21.73, 25.32),
    @(22, 22.89, 27.21),
    @(23, 21.68, 26.2),
    @(24, 19.47, 24.13),
    @(25, 23.88, 28.08),
    @(26, 23.66, 28.26),
    @(27, 37.25, 41.11),
    @(28, 14.89, 18.28)
)

foreach ($row in $data) {
    $r = $row[0]
    $fVal = $row[1]
    $gVal = $row[2]
    $ws.Cells.Item($r, 6).Value = $fVal
    $ws.Cells.Item($r, 7).Value = $gVal
}
